$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Hunk 0: sheet ALC
$ws_ALC.Cells.Item(51, 8).Value = 4371.4287
$ws_ALC.Cells.Item(51, 9).Value = 3150
$ws_ALC.Cells.Item(51, 10).Value = 6000
$ws_ALC.Cells.Item(51, 11).Value = 3150
$ws_ALC.Cells.Item(51, 12).Value = 6000
$ws_ALC.Cells.Item(51, 13).Value = -2666
$ws_ALC.Cells.Item(51, 14).Value = -6968

# Hunk 1: sheet ALC
$ws_ALC.Cells.Item(92, 8).Value = 849.71875
$ws_ALC.Cells.Item(92, 9).Value = 768.4167
$ws_ALC.Cells.Item(92, 10).Value = 1093.625
$ws_ALC.Cells.Item(92, 11).Value = 768.4167
$ws_ALC.Cells.Item(92, 12).Value = 1093.625
$ws_ALC.Cells.Item(92, 13).Value = 479.5833
$ws_ALC.Cells.Item(92, 14).Value = -3589.625

# Hunk 2: sheet ALC
$ws_ALC.Cells.Item(98, 8).Value = 5800
$ws_ALC.Cells.Item(98, 9).Value = 1000
$ws_ALC.Cells.Item(98, 10).Value = 9000
$ws_ALC.Cells.Item(98, 11).Value = 1000
$ws_ALC.Cells.Item(98, 12).Value = 9000
$ws_ALC.Cells.Item(98, 13).Value = 498
$ws_ALC.Cells.Item(98, 14).Value = -11996

# Hunk 3: sheet ALC
$ws_ALC.Cells.Item(122, 8).Value = 5800
$ws_ALC.Cells.Item(122, 9).Value = 1000
$ws_ALC.Cells.Item(122, 10).Value = 9000
$ws_ALC.Cells.Item(122, 11).Value = 3000
$ws_ALC.Cells.Item(122, 12).Value = 27000
$ws_ALC.Cells.Item(122, 13).Value = -550
$ws_ALC.Cells.Item(122, 14).Value = -31900

# Hunk 4: sheet ARM
$ws_ARM.Cells.Item(125, 8).Value = 35989.445
$ws_ARM.Cells.Item(125, 10).Value = 35989.445
$ws_ARM.Cells.Item(125, 12).Value = 35989.445
$ws_ARM.Cells.Item(125, 14).Value = -45829.445

# Hunk 5: sheet ARM
$ws_ARM.Cells.Item(132, 8).Value = 6907.838
$ws_ARM.Cells.Item(132, 9).Value = 4675.5864
$ws_ARM.Cells.Item(132, 10).Value = 14999.75
$ws_ARM.Cells.Item(132, 11).Value = 14026.7592
$ws_ARM.Cells.Item(132, 12).Value = 44999.25
$ws_ARM.Cells.Item(132, 13).Value = -11496.7592
$ws_ARM.Cells.Item(132, 14).Value = -50059.25

# Hunk 6: sheet BSM
$ws_BSM.Cells.Item(105, 8).Value = 62503384
$ws_BSM.Cells.Item(105, 9).Value = 76926160
$ws_BSM.Cells.Item(105, 11).Value = 76926160
$ws_BSM.Cells.Item(105, 13).Value = -76924413

# Hunk 7: sheet BSM
$ws_BSM.Cells.Item(123, 8).Value = 23982.857
$ws_BSM.Cells.Item(123, 10).Value = 23982.857
$ws_BSM.Cells.Item(123, 12).Value = 23982.857
$ws_BSM.Cells.Item(123, 14).Value = -33782.857

# Hunk 8: sheet CRP
$ws_CRP.Cells.Item(31, 8).Value = 18411.568
$ws_CRP.Cells.Item(31, 9).Value = 27086.857
$ws_CRP.Cells.Item(31, 10).Value = 2569.739
$ws_CRP.Cells.Item(31, 11).Value = 27086.857
$ws_CRP.Cells.Item(31, 12).Value = 2569.739
$ws_CRP.Cells.Item(31, 13).Value = -26791.857
$ws_CRP.Cells.Item(31, 14).Value = -3159.739

# Hunk 9: sheet CRP
$ws_CRP.Cells.Item(34, 8).Value = 18411.568
$ws_CRP.Cells.Item(34, 9).Value = 27086.857
$ws_CRP.Cells.Item(34, 10).Value = 2569.739
$ws_CRP.Cells.Item(34, 11).Value = 27086.857
$ws_CRP.Cells.Item(34, 12).Value = 2569.739
$ws_CRP.Cells.Item(34, 13).Value = -26884.857
$ws_CRP.Cells.Item(34, 14).Value = -2973.739

# Hunk 10: sheet CRP
$ws_CRP.Cells.Item(58, 8).Value = 1060.2653
$ws_CRP.Cells.Item(58, 9).Value = 1098.8462
$ws_CRP.Cells.Item(58, 10).Value = 909.8
$ws_CRP.Cells.Item(58, 11).Value = 1098.8462
$ws_CRP.Cells.Item(58, 12).Value = 909.8
$ws_CRP.Cells.Item(58, 13).Value = -895.8462
$ws_CRP.Cells.Item(58, 14).Value = -1315.8

# Hunk 11: sheet CRP
$ws_CRP.Cells.Item(64, 8).Value = 30271
$ws_CRP.Cells.Item(64, 10).Value = 30271
$ws_CRP.Cells.Item(64, 12).Value = 30271
$ws_CRP.Cells.Item(64, 14).Value = -30767

# Hunk 12: sheet CRP
$ws_CRP.Cells.Item(67, 8).Value = 30271
$ws_CRP.Cells.Item(67, 10).Value = 30271
$ws_CRP.Cells.Item(67, 12).Value = 30271
$ws_CRP.Cells.Item(67, 14).Value = -31987

# Hunk 13: sheet CRP
$ws_CRP.Cells.Item(68, 8).Value = 39813.637
$ws_CRP.Cells.Item(68, 10).Value = 39813.637
$ws_CRP.Cells.Item(68, 12).Value = 39813.637
$ws_CRP.Cells.Item(68, 14).Value = -41311.637

# Hunk 14: sheet CRP
$ws_CRP.Cells.Item(71, 8).Value = 39813.637
$ws_CRP.Cells.Item(71, 10).Value = 39813.637
$ws_CRP.Cells.Item(71, 12).Value = 119440.911
$ws_CRP.Cells.Item(71, 14).Value = -126928.911

# Hunk 15: sheet CRP
$ws_CRP.Cells.Item(81, 8).Value = 40328
$ws_CRP.Cells.Item(81, 10).Value = 40328
$ws_CRP.Cells.Item(81, 12).Value = 40328
$ws_CRP.Cells.Item(81, 14).Value = -42324

# Hunk 16: sheet CRP
$ws_CRP.Cells.Item(84, 8).Value = 40328
$ws_CRP.Cells.Item(84, 10).Value = 40328
$ws_CRP.Cells.Item(84, 12).Value = 120984
$ws_CRP.Cells.Item(84, 14).Value = -130968

# Hunk 17: sheet CRP
$ws_CRP.Cells.Item(122, 8).Value = 2187.2354
$ws_CRP.Cells.Item(122, 9).Value = 1935.2142
$ws_CRP.Cells.Item(122, 10).Value = 3363.3333
$ws_CRP.Cells.Item(122, 11).Value = 5805.642599999999
$ws_CRP.Cells.Item(122, 12).Value = 10089.9999
$ws_CRP.Cells.Item(122, 13).Value = -3355.642599999999
$ws_CRP.Cells.Item(122, 14).Value = -14989.9999

# Hunk 18: sheet CRP
$ws_CRP.Cells.Item(136, 8).Value = 1060.2653
$ws_CRP.Cells.Item(136, 9).Value = 1098.8462
$ws_CRP.Cells.Item(136, 10).Value = 909.8
$ws_CRP.Cells.Item(136, 11).Value = 3296.5386
$ws_CRP.Cells.Item(136, 12).Value = 2729.4
$ws_CRP.Cells.Item(136, 13).Value = -746.5385999999999
$ws_CRP.Cells.Item(136, 14).Value = -7829.4

# Hunk 19: sheet GSM
$ws_GSM.Cells.Item(70, 8).Value = 4715.485
$ws_GSM.Cells.Item(70, 9).Value = 4351.143
$ws_GSM.Cells.Item(70, 10).Value = 4983.9473
$ws_GSM.Cells.Item(70, 11).Value = 4351.143
$ws_GSM.Cells.Item(70, 12).Value = 4983.9473
$ws_GSM.Cells.Item(70, 13).Value = -4081.143
$ws_GSM.Cells.Item(70, 14).Value = -5523.9473

# Hunk 20: sheet GSM
$ws_GSM.Cells.Item(73, 8).Value = 4715.485
$ws_GSM.Cells.Item(73, 9).Value = 4351.143
$ws_GSM.Cells.Item(73, 10).Value = 4983.9473
$ws_GSM.Cells.Item(73, 11).Value = 4351.143
$ws_GSM.Cells.Item(73, 12).Value = 4983.9473
$ws_GSM.Cells.Item(73, 13).Value = -3415.143
$ws_GSM.Cells.Item(73, 14).Value = -6855.9473

# Hunk 21: sheet GSM
$ws_GSM.Cells.Item(102, 8).Value = 1367.826
$ws_GSM.Cells.Item(102, 9).Value = 1154.1765
$ws_GSM.Cells.Item(102, 10).Value = 1973.1666
$ws_GSM.Cells.Item(102, 11).Value = 1154.1765
$ws_GSM.Cells.Item(102, 12).Value = 1973.1666
$ws_GSM.Cells.Item(102, 13).Value = 467.8235
$ws_GSM.Cells.Item(102, 14).Value = -5217.1666

# Hunk 22: sheet GSM
$ws_GSM.Cells.Item(122, 8).Value = 2180.5715
$ws_GSM.Cells.Item(122, 9).Value = 2052.8
$ws_GSM.Cells.Item(122, 10).Value = 2500
$ws_GSM.Cells.Item(122, 11).Value = 6158.400000000001
$ws_GSM.Cells.Item(122, 12).Value = 7500
$ws_GSM.Cells.Item(122, 13).Value = -3708.400000000001
$ws_GSM.Cells.Item(122, 14).Value = -12400

# Hunk 23: sheet GSM
$ws_GSM.Cells.Item(126, 8).Value = 1305.7142
$ws_GSM.Cells.Item(126, 9).Value = 1236
$ws_GSM.Cells.Item(126, 11).Value = 3708
$ws_GSM.Cells.Item(126, 13).Value = -1238

# Hunk 24: sheet GSM
$ws_GSM.Cells.Item(132, 8).Value = 6061.102
$ws_GSM.Cells.Item(132, 9).Value = 4385.1353
$ws_GSM.Cells.Item(132, 10).Value = 11228.667
$ws_GSM.Cells.Item(132, 11).Value = 13155.4059
$ws_GSM.Cells.Item(132, 12).Value = 33686.001
$ws_GSM.Cells.Item(132, 13).Value = -10625.4059
$ws_GSM.Cells.Item(132, 14).Value = -38746.001

# Hunk 25: sheet LTW
$ws_LTW.Cells.Item(40, 8).Value = 5505.5
$ws_LTW.Cells.Item(40, 9).Value = 4749.9
$ws_LTW.Cells.Item(40, 11).Value = 4749.9
$ws_LTW.Cells.Item(40, 13).Value = -4613.9

# Hunk 26: sheet LTW
$ws_LTW.Cells.Item(64, 8).Value = 18690
$ws_LTW.Cells.Item(64, 10).Value = 18690
$ws_LTW.Cells.Item(64, 12).Value = 18690
$ws_LTW.Cells.Item(64, 14).Value = -19140

# Hunk 27: sheet LTW
$ws_LTW.Cells.Item(67, 8).Value = 18690
$ws_LTW.Cells.Item(67, 10).Value = 18690
$ws_LTW.Cells.Item(67, 12).Value = 18690
$ws_LTW.Cells.Item(67, 14).Value = -20250

# Hunk 28: sheet LTW
$ws_LTW.Cells.Item(120, 8).Value = 38000
$ws_LTW.Cells.Item(120, 10).Value = 38000
$ws_LTW.Cells.Item(120, 12).Value = 38000
$ws_LTW.Cells.Item(120, 14).Value = -47676

# Hunk 29: sheet LTW
$ws_LTW.Cells.Item(122, 8).Value = 3814.2903
$ws_LTW.Cells.Item(122, 9).Value = 4081.8096
$ws_LTW.Cells.Item(122, 11).Value = 12245.4288
$ws_LTW.Cells.Item(122, 13).Value = -9795.4288

# Hunk 30: sheet LTW
$ws_LTW.Cells.Item(132, 8).Value = 4202.448
$ws_LTW.Cells.Item(132, 9).Value = 4391.5684
$ws_LTW.Cells.Item(132, 10).Value = 3608.0715
$ws_LTW.Cells.Item(132, 11).Value = 13174.7052
$ws_LTW.Cells.Item(132, 12).Value = 10824.2145
$ws_LTW.Cells.Item(132, 13).Value = -10644.7052
$ws_LTW.Cells.Item(132, 14).Value = -15884.2145

# Hunk 31: sheet WVR
$ws_WVR.Cells.Item(63, 8).Value = 7035.2856
$ws_WVR.Cells.Item(63, 10).Value = 7035.2856
$ws_WVR.Cells.Item(63, 12).Value = 7035.2856
$ws_WVR.Cells.Item(63, 14).Value = -8283.285599999999

# Hunk 32: sheet WVR
$ws_WVR.Cells.Item(66, 8).Value = 7035.2856
$ws_WVR.Cells.Item(66, 10).Value = 7035.2856
$ws_WVR.Cells.Item(66, 12).Value = 21105.8568
$ws_WVR.Cells.Item(66, 14).Value = -27345.8568

# Hunk 33: sheet WVR
$ws_WVR.Cells.Item(107, 8).Value = 495
$ws_WVR.Cells.Item(107, 9).Value = 475.4
$ws_WVR.Cells.Item(107, 10).Value = 509
$ws_WVR.Cells.Item(107, 11).Value = 1426.2
$ws_WVR.Cells.Item(107, 12).Value = 1527
$ws_WVR.Cells.Item(107, 13).Value = 493.8000000000002
$ws_WVR.Cells.Item(107, 14).Value = -5367

# Hunk 34: sheet WVR
$ws_WVR.Cells.Item(122, 8).Value = 3049679.8
$ws_WVR.Cells.Item(122, 9).Value = 2116765.8
$ws_WVR.Cells.Item(122, 10).Value = 4169177
$ws_WVR.Cells.Item(122, 11).Value = 6350297.399999999
$ws_WVR.Cells.Item(122, 12).Value = 12507531
$ws_WVR.Cells.Item(122, 13).Value = -6347847.399999999
$ws_WVR.Cells.Item(122, 14).Value = -12512431

# Hunk 35: sheet WVR
$ws_WVR.Cells.Item(132, 8).Value = 6419.15
$ws_WVR.Cells.Item(132, 9).Value = 7470.7188
$ws_WVR.Cells.Item(132, 10).Value = 2212.875
$ws_WVR.Cells.Item(132, 11).Value = 22412.1564
$ws_WVR.Cells.Item(132, 12).Value = 6638.625
$ws_WVR.Cells.Item(132, 13).Value = -19882.1564
$ws_WVR.Cells.Item(132, 14).Value = -11698.625
